# Adds a new weekly price record for "Albahaca" (Terminal La Palmera de
# La Serena) dated 2021-11-16 (serial 44516). The new record is inserted
# as row 56, pushing the previous row 56 (dated 2021-11-09 / serial 44509)
# down to row 57, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 56; this shifts the former row 56
# (and its formatting) down to row 57.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly record.
$ws.Cells.Item(56, 1).Value = 8
$ws.Cells.Item(56, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(56, 3).Value = "Coquimbo"
$ws.Cells.Item(56, 4).Value = 44516
$ws.Cells.Item(56, 5).Value = 4
$ws.Cells.Item(56, 6).Value = 100112052
$ws.Cells.Item(56, 7).Value = "Albahaca"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 740
$ws.Cells.Item(56, 11).Value = 3000
$ws.Cells.Item(56, 12).Value = 4000
$ws.Cells.Item(56, 13).Value = 3500
$ws.Cells.Item(56, 14).Value = "`$/paquete"
$ws.Cells.Item(56, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(56, 16).Value = 3500
$ws.Cells.Item(56, 17).Value = 1
$ws.Cells.Item(56, 18).Value = "Hortaliza"
